$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from existing adjacent column J to new column K
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)  # xlPasteFormats

# Add new column K (year 2023) matching the existing table pattern
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1533.3
$ws.Range("K5").Value = 944.8
$ws.Range("K6").Value = 1914.8
